# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#   - Status cells flip from "Ready for handoff" to "Handed back: in sync with en-US"
#   - Each language sheet's "Latest Target File" (I) / "Latest Handback File" (J) /
#     "Latest Handback DateTime" (K) columns get filled in for the two data rows
#   - A hyperlink (matching the one already on column A) is added on the new
#     "Latest Target File" links
#   - A few columns are widened to fit the newly-populated long text

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdDisplay  = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md"
$mdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25ddfdab2bb48d46b7ef5852d591145e567e9598/e2e/b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md"

# --- Overview sheet: Status columns (E/F) for both data rows -------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Widen the two Status columns to fit the longer text
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet -----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = $mdDisplay
$wsZh.Range("I3").Value = $mdDisplay
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay)

$wsZh.Range("J2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.zh-cn.xlf"
$wsZh.Range("J3").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-19 01:00:31"
$wsZh.Range("K3").Value = "2016-08-19 01:00:31"

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de sheet -----------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = $mdDisplay
$wsDe.Range("I3").Value = $mdDisplay
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay)

$wsDe.Range("J2").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.de-de.xlf"
$wsDe.Range("J3").Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-19 01:00:39"
$wsDe.Range("K3").Value = "2016-08-19 01:00:39"

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
